$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.986.16"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "3.390.09"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'559.90"
$ws.Range("E5").Value = "  +3.56%  "
$ws.Range("D6").Value = "'175.16"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("E7").Value = "  +2.66%  "
$ws.Range("D8").Value = "3.377.88"
$ws.Range("E8").Value = "  +2.79%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +9.36%  "
$ws.Range("D11").Value = "'0.632"
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("D12").Value = "'54.21"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("E13").Value = "  +5.94%  "
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").Value = "3.932.93"
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "3.378.33"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").Value = "64.877.58"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").Value = "'11.81"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("E21").Value = "  +3.26%  "
$ws.Range("D22").Value = "'465.29"
$ws.Range("E22").Value = "  +12.06%  "
$ws.Range("D23").Value = "'4.89"
$ws.Range("E23").Value = "  +11.60%  "
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").Value = "'86.49"
$ws.Range("E25").Value = "  +4.80%  "
$ws.Range("D26").Value = "'13.54"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").Value = "'2.96"
$ws.Range("E27").Value = "  +9.39%  "
$ws.Range("D28").Value = "'10.84"
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("E29").Value = "  +2.62%  "
$ws.Range("D30").Value = "'30.83"
$ws.Range("E30").Value = "  +6.75%  "
$ws.Range("D31").Value = "'6.76"
$ws.Range("E31").Value = "  +6.97%  "
$ws.Range("D32").Value = "'11.49"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Value = "'572.00"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'61.29"
$ws.Range("E34").Value = "  +6.09%  "
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'3.64"
$ws.Range("E37").Value = "  +7.24%  "
$ws.Range("D38").Value = "'0.140"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").Value = "0.0₃0746"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").Value = "'0.370"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").Value = "3.088.02"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +4.04%  "
$ws.Range("D45").Value = "'0.0416"
$ws.Range("E45").Value = "  +4.68%  "
$ws.Range("E46").Value = "  +5.59%  "
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").Value = "'3.13"
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "'138.19"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("D51").Value = "'8.24"
$ws.Range("E51").Value = "  +3.28%  "
